# Daily attendance processing - 2026-01-20 08:06:48
# The "Recorded By" column (G) lists session recorders such as
# "System, dnasr281@gmail.com". Swap the order of the two names so the
# individual's email is listed first, e.g. "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

$recordedByRange = $ws.Range("G1:G$lastRow")
[void]$recordedByRange.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System", -1, 1, $false, $false, $true, $true)
